# Generate Report for handoff
# Update the "Latest Handoff Datetime" column (D) for the d1a32dc4 row (row 5)
# on both the zh-cn and de-de sheets, recording the datetime the report
# was (re)generated / handed off.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D5").Value = "2016-01-22 02:07:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D5").Value = "2016-01-22 02:07:48"
